$d = $word.ActiveDocument

$pairs = @(
    @("160×6=", "392×3="),
    @("342×9=", "568×7="),
    @("415×8=", "715×8="),
    @("856×4=", "367×8="),
    @("731×4=", "175×5="),
    @("424×2=", "163×6="),
    @("111×4=", "195×6="),
    @("863×5=", "558×9="),
    @("110×6=", "732×2="),
    @("992×9=", "812×4="),
    @("245×2=", "106×3="),
    @("447×5=", "262×4="),
    @("980×3=", "906×3="),
    @("750×7=", "545×2="),
    @("726×8=", "519×6="),
    @("349×9=", "389×4="),
    @("118×9=", "190×4="),
    @("579×2=", "619×3="),
    @("481×9=", "296×5="),
    @("875×8=", "491×2="),
    @("983×2=", "672×4="),
    @("269×9=", "844×8="),
    @("423×2=", "429×5="),
    @("509×4=", "623×9="),
    @("437×4=", "537×6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
